$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.523.32"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "1.637.93"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'307.28"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.3768"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "'0.3651"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "'1.271"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "'0.08181"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'23.03"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'6.645"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'0.00001279"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D17").Value = "1.638.50"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "'94.77"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "'0.06942"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "'6.559"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "23.522.77"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'12.82"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "'3.098"
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("D26").Value = "'2.424"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'151.62"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'5.360"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "'135.59"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").Value = "'2.367"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "1.820.89"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "'0.9656"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'0.02822"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").Value = "'10.34"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'0.07363"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").Value = "'0.2542"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").Value = "'6.192"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "'0.08870"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "'1.382"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "'0.7113"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'16.19"
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("D45").Value = "'0.6549"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'2.345"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'4.040"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "'0.07977"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'129.45"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").Value = "'1.209"
$ws.Range("E51").Value = "  +0.49%  "
